$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 153, shifting the existing
# rows 153-159 down to 155-161.
$ws.Rows.Item(153).Insert()
$ws.Rows.Item(153).Insert()

# New row 153: "Extra" quality entry dated 2022-01-17 (serial 44578)
$ws.Cells.Item(153, 1).Value = 5
$ws.Cells.Item(153, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(153, 3).Value = "Maule"
$ws.Cells.Item(153, 4).Value = 44578
$ws.Cells.Item(153, 5).Value = 7
$ws.Cells.Item(153, 6).Value = 100112028
$ws.Cells.Item(153, 7).Value = "Sandia"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Extra"
$ws.Cells.Item(153, 10).Value = 2500
$ws.Cells.Item(153, 11).Value = 2000
$ws.Cells.Item(153, 12).Value = 2000
$ws.Cells.Item(153, 13).Value = 2000
$ws.Cells.Item(153, 14).Value = "$/unidad"
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 2000
$ws.Cells.Item(153, 17).Value = 1
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# New row 154: "Primera" quality entry dated 2022-01-17 (serial 44578)
$ws.Cells.Item(154, 1).Value = 5
$ws.Cells.Item(154, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(154, 3).Value = "Maule"
$ws.Cells.Item(154, 4).Value = 44578
$ws.Cells.Item(154, 5).Value = 7
$ws.Cells.Item(154, 6).Value = 100112028
$ws.Cells.Item(154, 7).Value = "Sandia"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 2000
$ws.Cells.Item(154, 11).Value = 1500
$ws.Cells.Item(154, 12).Value = 1500
$ws.Cells.Item(154, 13).Value = 1500
$ws.Cells.Item(154, 14).Value = "$/unidad"
$ws.Cells.Item(154, 15).Value = "Región del Maule"
$ws.Cells.Item(154, 16).Value = 1500
$ws.Cells.Item(154, 17).Value = 1
$ws.Cells.Item(154, 18).Value = "Hortaliza"
